$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap the three pairs of rows that got re-ordered upstream ---
# (id column A keeps its original sequential numbering; every other
# column's data moves between the two rows in the pair)

$swapCols = @("B","F","G","H","I","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Swap-RowData($r1, $r2) {
    foreach ($col in $swapCols) {
        $c1 = $ws.Range($col + $r1)
        $c2 = $ws.Range($col + $r2)
        $v1 = $c1.Value()
        $v2 = $c2.Value()
        $c1.Value = $v2
        $c2.Value = $v1
    }
}

Swap-RowData 147 149
Swap-RowData 150 152
Swap-RowData 153 154

# --- Step 2: append 9 new fixtures (rows 323-331) ---

$ws.Range("A323").Value = 321
$ws.Range("B323").Value = 6774883
$ws.Range("C323").Value = "Poland Ekstraklasa"
$ws.Range("D323").Value = "Poland Ekstraklasa"
$ws.Range("E323").Value = 45331.58333333334
$ws.Range("F323").Value = "Stal Mielec"
$ws.Range("G323").Value = "MKS Puszcza Niepolomice"
$ws.Range("K323").Value = 2
$ws.Range("L323").Value = 3.5
$ws.Range("M323").Value = 3.6
$ws.Range("N323").Value = 2.1
$ws.Range("O323").Value = 3.4
$ws.Range("P323").Value = 3.3
$ws.Range("Q323").Value = -0.25
$ws.Range("R323").Value = 1.875
$ws.Range("S323").Value = 1.975
$ws.Range("T323").Value = 2.5
$ws.Range("U323").Value = 2.05
$ws.Range("V323").Value = 1.8
$ws.Range("W323").Value = 0
$ws.Range("X323").Value = 0
$ws.Range("Y323").Value = 0
$ws.Range("Z323").Value = 0
$ws.Range("AA323").Value = 0

$ws.Range("A324").Value = 322
$ws.Range("B324").Value = 6774458
$ws.Range("C324").Value = "Poland Ekstraklasa"
$ws.Range("D324").Value = "Poland Ekstraklasa"
$ws.Range("E324").Value = 45331.6875
$ws.Range("F324").Value = "Ruch Chorzow"
$ws.Range("G324").Value = "Legia Warsaw"
$ws.Range("K324").Value = 5
$ws.Range("L324").Value = 3.75
$ws.Range("M324").Value = 1.666
$ws.Range("N324").Value = 5.5
$ws.Range("O324").Value = 3.8
$ws.Range("P324").Value = 1.615
$ws.Range("Q324").Value = 0.75
$ws.Range("R324").Value = 2.05
$ws.Range("S324").Value = 1.8
$ws.Range("T324").Value = 2.5
$ws.Range("U324").Value = 1.9
$ws.Range("V324").Value = 1.95
$ws.Range("W324").Value = 0
$ws.Range("X324").Value = 0
$ws.Range("Y324").Value = 0
$ws.Range("Z324").Value = 0
$ws.Range("AA324").Value = 0

$ws.Range("A325").Value = 323
$ws.Range("B325").Value = 6775535
$ws.Range("C325").Value = "Poland Ekstraklasa"
$ws.Range("D325").Value = "Poland Ekstraklasa"
$ws.Range("E325").Value = 45332.45833333334
$ws.Range("F325").Value = "Cracovia Krakow"
$ws.Range("G325").Value = "Radomiak Radom"
$ws.Range("K325").Value = 2.15
$ws.Range("L325").Value = 3.3
$ws.Range("M325").Value = 3.4
$ws.Range("N325").Value = 2.15
$ws.Range("O325").Value = 3.3
$ws.Range("P325").Value = 3.4
$ws.Range("Q325").Value = -0.25
$ws.Range("R325").Value = 1.875
$ws.Range("S325").Value = 1.975
$ws.Range("T325").Value = 2.25
$ws.Range("U325").Value = 2
$ws.Range("V325").Value = 1.85
$ws.Range("W325").Value = 0
$ws.Range("X325").Value = 0
$ws.Range("Y325").Value = 0
$ws.Range("Z325").Value = 0
$ws.Range("AA325").Value = 0

$ws.Range("A326").Value = 324
$ws.Range("B326").Value = 6775537
$ws.Range("C326").Value = "Poland Ekstraklasa"
$ws.Range("D326").Value = "Poland Ekstraklasa"
$ws.Range("E326").Value = 45332.5625
$ws.Range("F326").Value = "Piast Gliwice"
$ws.Range("G326").Value = "Gornik Zabrze"
$ws.Range("K326").Value = 1.909
$ws.Range("L326").Value = 3.3
$ws.Range("M326").Value = 4.2
$ws.Range("N326").Value = 1.909
$ws.Range("O326").Value = 3.3
$ws.Range("P326").Value = 4.2
$ws.Range("Q326").Value = -0.5
$ws.Range("R326").Value = 1.9
$ws.Range("S326").Value = 1.95
$ws.Range("T326").Value = 2.25
$ws.Range("U326").Value = 2.05
$ws.Range("V326").Value = 1.8
$ws.Range("W326").Value = 0
$ws.Range("X326").Value = 0
$ws.Range("Y326").Value = 0
$ws.Range("Z326").Value = 0
$ws.Range("AA326").Value = 0

$ws.Range("A327").Value = 325
$ws.Range("B327").Value = 6775536
$ws.Range("C327").Value = "Poland Ekstraklasa"
$ws.Range("D327").Value = "Poland Ekstraklasa"
$ws.Range("E327").Value = 45332.66666666666
$ws.Range("F327").Value = "Lech Poznan"
$ws.Range("G327").Value = "Zaglebie Lubin"
$ws.Range("K327").Value = 1.571
$ws.Range("L327").Value = 4
$ws.Range("M327").Value = 5.5
$ws.Range("N327").Value = 1.571
$ws.Range("O327").Value = 4
$ws.Range("P327").Value = 5.5
$ws.Range("Q327").Value = -1
$ws.Range("R327").Value = 2.025
$ws.Range("S327").Value = 1.825
$ws.Range("T327").Value = 2.5
$ws.Range("U327").Value = 1.85
$ws.Range("V327").Value = 2
$ws.Range("W327").Value = 0
$ws.Range("X327").Value = 0
$ws.Range("Y327").Value = 0
$ws.Range("Z327").Value = 0
$ws.Range("AA327").Value = 0

$ws.Range("A328").Value = 326
$ws.Range("B328").Value = 6775539
$ws.Range("C328").Value = "Poland Ekstraklasa"
$ws.Range("D328").Value = "Poland Ekstraklasa"
$ws.Range("E328").Value = 45333.35416666666
$ws.Range("F328").Value = "Warta Poznan"
$ws.Range("G328").Value = "Rakow Czestochowa"
$ws.Range("K328").Value = 4.75
$ws.Range("L328").Value = 3.5
$ws.Range("M328").Value = 1.75
$ws.Range("N328").Value = 4.75
$ws.Range("O328").Value = 3.5
$ws.Range("P328").Value = 1.75
$ws.Range("Q328").Value = 0.75
$ws.Range("R328").Value = 1.825
$ws.Range("S328").Value = 2.025
$ws.Range("T328").Value = 2.25
$ws.Range("U328").Value = 1.925
$ws.Range("V328").Value = 1.925
$ws.Range("W328").Value = 0
$ws.Range("X328").Value = 0
$ws.Range("Y328").Value = 0
$ws.Range("Z328").Value = 0
$ws.Range("AA328").Value = 0

$ws.Range("A329").Value = 327
$ws.Range("B329").Value = 6775540
$ws.Range("C329").Value = "Poland Ekstraklasa"
$ws.Range("D329").Value = "Poland Ekstraklasa"
$ws.Range("E329").Value = 45333.45833333334
$ws.Range("F329").Value = "Widzew Lodz"
$ws.Range("G329").Value = "Jagiellonia Bialystok"
$ws.Range("K329").Value = 2.6
$ws.Range("L329").Value = 3.4
$ws.Range("M329").Value = 2.6
$ws.Range("N329").Value = 2.7
$ws.Range("O329").Value = 3.4
$ws.Range("P329").Value = 2.5
$ws.Range("Q329").Value = 0
$ws.Range("R329").Value = 2.025
$ws.Range("S329").Value = 1.825
$ws.Range("T329").Value = 2.5
$ws.Range("U329").Value = 1.925
$ws.Range("V329").Value = 1.925
$ws.Range("W329").Value = 0
$ws.Range("X329").Value = 0
$ws.Range("Y329").Value = 0
$ws.Range("Z329").Value = 0
$ws.Range("AA329").Value = 0

$ws.Range("A330").Value = 328
$ws.Range("B330").Value = 6775538
$ws.Range("C330").Value = "Poland Ekstraklasa"
$ws.Range("D330").Value = "Poland Ekstraklasa"
$ws.Range("E330").Value = 45333.5625
$ws.Range("F330").Value = "Slask Wroclaw"
$ws.Range("G330").Value = "Pogon Szczecin"
$ws.Range("K330").Value = 3.9
$ws.Range("L330").Value = 3.5
$ws.Range("M330").Value = 1.909
$ws.Range("N330").Value = 2.75
$ws.Range("O330").Value = 3.3
$ws.Range("P330").Value = 2.55
$ws.Range("Q330").Value = 0
$ws.Range("R330").Value = 2
$ws.Range("S330").Value = 1.85
$ws.Range("T330").Value = 2.5
$ws.Range("U330").Value = 1.975
$ws.Range("V330").Value = 1.875
$ws.Range("W330").Value = 0
$ws.Range("X330").Value = 0
$ws.Range("Y330").Value = 0
$ws.Range("Z330").Value = 0
$ws.Range("AA330").Value = 0

$ws.Range("A331").Value = 329
$ws.Range("B331").Value = 6774457
$ws.Range("C331").Value = "Poland Ekstraklasa"
$ws.Range("D331").Value = "Poland Ekstraklasa"
$ws.Range("E331").Value = 45334.625
$ws.Range("F331").Value = "Korona Kielce"
$ws.Range("G331").Value = "LKS Lodz"
$ws.Range("K331").Value = 1.95
$ws.Range("L331").Value = 3.5
$ws.Range("M331").Value = 3.75
$ws.Range("N331").Value = 1.75
$ws.Range("O331").Value = 3.4
$ws.Range("P331").Value = 5
$ws.Range("Q331").Value = -0.75
$ws.Range("R331").Value = 2.025
$ws.Range("S331").Value = 1.825
$ws.Range("T331").Value = 2.5
$ws.Range("U331").Value = 2.025
$ws.Range("V331").Value = 1.825
$ws.Range("W331").Value = 0
$ws.Range("X331").Value = 0
$ws.Range("Y331").Value = 0
$ws.Range("Z331").Value = 0
$ws.Range("AA331").Value = 0


# --- Step 3: apply the correct cell styles to the new rows (A -> bold/border
# style used for the id column, E -> date/time number format) by copying
# formats from the previous last row (322), which already carries them. ---
$ws.Range("A322").Copy() | Out-Null
$ws.Range("A323:A331").PasteSpecial(-4122) | Out-Null
$ws.Range("E322").Copy() | Out-Null
$ws.Range("E323:E331").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

"done"
